{"js": "// transaction.delete.humanResource.setBusinessTripCostComponent -> setWorkTimeEpoch\n// (API renamed, description/date/sample-payload text refreshed accordingly)\n\nasync function replaceAll(searchText, replacementText) {\n  const results = context.document.body.search(searchText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(replacementText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n  return results.items.length;\n}\n\n// 1. API title (\"transaction.delete.humanResource. setBusinessTripCostComponent\"\n//    -> \"transaction.delete.humanResource.setWorkTimeEpoch\") \u2014 handle the lone\n//    occurrence that has a space before \"set\" first so the space is dropped.\nawait replaceAll(\". setBusinessTripCostComponent\", \".setWorkTimeEpoch\");\n\n// 2. Every remaining \"BusinessTripCostComponent\" (directory paths, file name,\n//    JSON schema text, sample code, etc.) simply becomes \"WorkTimeEpoch\".\nawait replaceAll(\"BusinessTripCostComponent\", \"WorkTimeEpoch\");\n\n// 3. Indonesian description of the endpoint.\nawait replaceAll(\n  \"Menghapus Data Komponen Biaya Perjalanan Bisnis\",\n  \"Menghapus Data Masa Waktu Kerja\"\n);\n\n// 4. Document date bumped from the 23rd to the 25th.\nawait replaceAll(\"23\", \"25\");\n\n// 5. Sample API web token refreshed (new \"iat\" + signature).\nawait replaceAll(\n  \"eyJhbGciOiJIUzI1NiIsInR5cCI6IkpXVCJ9.eyJsb2dnZWRJbkFzIjoidGVndWgucHJhdGFtYSIsImlhdCI6MTYwNjA5Nzg4MH0.d1AB_XF31WOFS7dhxvEHyJmPybR5ju4YHiuF_ZbSf5Q\",\n  \"eyJhbGciOiJIUzI1NiIsInR5cCI6IkpXVCJ9.eyJsb2dnZWRJbkFzIjoidGVndWgucHJhdGFtYSIsImlhdCI6MTYwNjI2OTA1NH0.NjJJegg6WRVQ3LHksbKcni92MkyzjfYpxzrFvgLu2FQ\"\n);\n\n// 6. Sample recordID value changed from 81... to 33...\nawait replaceAll(\"81\", \"33\");\n\n// 7. Word moves its \"last edit\" (_GoBack) bookmark to the spot that was just\n//    typed over; replicate that by re-pointing it at the new \"33\" run.\ntry {\n  context.document.deleteBookmark(\"_GoBack\");\n  await context.sync();\n} catch (e) {\n  // no-op: nothing to delete\n}\n\nconst recordIdHits = context.document.body.search(\"33\", { matchCase: true });\nrecordIdHits.load(\"items\");\nawait context.sync();\nif (recordIdHits.items.length > 0) {\n  const afterRange = recordIdHits.items[recordIdHits.items.length - 1].getRange(Word.RangeLocation.after);\n  afterRange.insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n", "ps1": "# transaction.delete.humanResource.setBusinessTripCostComponent -> setWorkTimeEpoch\n# (API renamed, description/date/sample-payload text refreshed accordingly)\n\n$d = $word.ActiveDocument\n\nfunction Replace-AllText($findText, $replaceText) {\n    $f = $d.Content.Find\n    $f.ClearFormatting()\n    $f.Replacement.ClearFormatting()\n    $f.Text = $findText\n    $f.Replacement.Text = $replaceText\n    # wdFindContinue = 1, wdReplaceAll = 2\n    $f.Execute($f.Text, $false, $false, $false, $false, $false, $true, 1, $false, $f.Replacement.Text, 2)\n}\n\n# 1. API title (\"transaction.delete.humanResource. setBusinessTripCostComponent\"\n#    -> \"transaction.delete.humanResource.setWorkTimeEpoch\") \u2014 handle the lone\n#    occurrence that has a space before \"set\" first so the space is dropped.\nReplace-AllText \". setBusinessTripCostComponent\" \".setWorkTimeEpoch\"\n\n# 2. Every remaining \"BusinessTripCostComponent\" (directory paths, file name,\n#    JSON schema text, sample code, etc.) simply becomes \"WorkTimeEpoch\".\nReplace-AllText \"BusinessTripCostComponent\" \"WorkTimeEpoch\"\n\n# 3. Indonesian description of the endpoint.\nReplace-AllText \"Menghapus Data Komponen Biaya Perjalanan Bisnis\" \"Menghapus Data Masa Waktu Kerja\"\n\n# 4. Document date bumped from the 23rd to the 25th.\nReplace-AllText \"23\" \"25\"\n\n# 5. Sample API web token refreshed (new \"iat\" + signature).\nReplace-AllText \"eyJhbGciOiJIUzI1NiIsInR5cCI6IkpXVCJ9.eyJsb2dnZWRJbkFzIjoidGVndWgucHJhdGFtYSIsImlhdCI6MTYwNjA5Nzg4MH0.d1AB_XF31WOFS7dhxvEHyJmPybR5ju4YHiuF_ZbSf5Q\" \"eyJhbGciOiJIUzI1NiIsInR5cCI6IkpXVCJ9.eyJsb2dnZWRJbkFzIjoidGVndWgucHJhdGFtYSIsImlhdCI6MTYwNjI2OTA1NH0.NjJJegg6WRVQ3LHksbKcni92MkyzjfYpxzrFvgLu2FQ\"\n\n# 6. Sample recordID value changed from 81... to 33...\nReplace-AllText \"81\" \"33\"\n\n# 7. Word moves its \"last edit\" (_GoBack) bookmark to the spot that was just\n#    typed over; replicate that by re-pointing it at the new \"33\" run.\n$goBack = $d.Bookmarks(\"_GoBack\")\n$goBack.Delete()\n\n$rng = $d.Content\n$rngFind = $rng.Find\n$rngFind.ClearFormatting()\n$rngFind.Text = \"33\"\n$rngFind.Execute() | Out-Null\n\n$bmRange = $d.Range($rng.End, $rng.End)\n$d.Bookmarks.Add(\"_GoBack\", $bmRange)\n"}
